# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder tied countries (same "Casos totales") that swapped position ---
# Belice <-> Nueva Caledonia (rows 193/194)
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

$ws.Range("A194").Value = "Belice"
$ws.Range("B194").Value = 18
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 16
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2

# Curazao <-> Dominica (rows 198/199)
$ws.Range("A198").Value = "Dominica"
$ws.Range("B198").Value = 16
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 15
$ws.Range("E198").Value = 1
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("B199").Value = 16
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 14
$ws.Range("E199").Value = 1
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

# Sahara Occidental <-> San Bartolome (rows 215/216) - data is identical, only label order swaps
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Sahara Occidental"

# --- Update Reunion's daily figures (row 128) ---
$ws.Range("B128").Value = 437
$ws.Range("C128").Value = 1
$ws.Range("E128").Value = 83
$ws.Range("F128").Value = 4

# --- Update the "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 02:35"
